# Update the LR-pairs sheet with the refreshed TPM-based values.
# Rows 2-5: sending cluster changes from "ECs" to "MuSCs" (same Proc/Thbd ligand-
# receptor pair), with refreshed expression/specificity numbers.
# Rows 6-9: newly added sending cluster "Resolving-Mac" x Proc/Thbd, one row per
# target cluster (ECs, FAPs, MuSCs, Resolving-Mac).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns: A Sending cluster, B Ligand symbol, C Receptor symbol, D Target cluster,
# E..T the associated expression / specificity metrics.
$rows = @(
    @("MuSCs", "Proc", "Thbd", "ECs", 2.0, 0.6666666666666666, 0.132605, 0.397815, 0.776265286171174, 0.776265286171174, 3.0, 1.0, 41.60833666666667, 124.82501, 0.5886423873735626, 0.5886423873735626, 5.517473483683333, 49.65726135315, 0.4569426512870216, 0.4569426512870216),
    @("MuSCs", "Proc", "Thbd", "FAPs", 2.0, 0.6666666666666666, 0.132605, 0.397815, 0.776265286171174, 0.776265286171174, 3.0, 1.0, 17.71526766666667, 53.145803, 0.2506218293658061, 0.2506218293658061, 2.349133068938333, 21.142197620445, 0.1945490260933906, 0.1945490260933906),
    @("MuSCs", "Proc", "Thbd", "MuSCs", 2.0, 0.6666666666666666, 0.132605, 0.397815, 0.776265286171174, 0.776265286171174, 3.0, 1.0, 6.935318333333334, 20.805955, 0.09811549001908282, 0.09811549001908279, 0.9196578875916668, 8.276920988325, 0.07616364893748828, 0.07616364893748827),
    @("MuSCs", "Proc", "Thbd", "Resolving-Mac", 2.0, 0.6666666666666666, 0.132605, 0.397815, 0.776265286171174, 0.776265286171174, 3.0, 1.0, 4.426331333333334, 13.278994, 0.06262029324154843, 0.06262029324154841, 0.5869536664566667, 5.28258299811, 0.04860995985327342, 0.04860995985327341),
    @("Resolving-Mac", "Proc", "Thbd", "ECs", 1.0, 0.3333333333333333, 0.03821933333333333, 0.114658, 0.2237347138288261, 0.2237347138288261, 3.0, 1.0, 41.60833666666667, 124.82501, 0.5886423873735626, 0.5886423873735626, 1.590242888508889, 14.31218599658, 0.131699736086541, 0.131699736086541),
    @("Resolving-Mac", "Proc", "Thbd", "FAPs", 1.0, 0.3333333333333333, 0.03821933333333333, 0.114658, 0.2237347138288261, 0.2237347138288261, 3.0, 1.0, 17.71526766666667, 53.145803, 0.2506218293658061, 0.2506218293658061, 0.6770657200415555, 6.093591480374, 0.05607280327241553, 0.05607280327241553),
    @("Resolving-Mac", "Proc", "Thbd", "MuSCs", 1.0, 0.3333333333333333, 0.03821933333333333, 0.114658, 0.2237347138288261, 0.2237347138288261, 3.0, 1.0, 6.935318333333334, 20.805955, 0.09811549001908282, 0.09811549001908279, 0.2650632431544445, 2.38556918839, 0.02195184108159454, 0.02195184108159453),
    @("Resolving-Mac", "Proc", "Thbd", "Resolving-Mac", 1.0, 0.3333333333333333, 0.03821933333333333, 0.114658, 0.2237347138288261, 0.2237347138288261, 3.0, 1.0, 4.426331333333334, 13.278994, 0.06262029324154843, 0.06262029324154841, 0.1691714326724445, 1.522542894052, 0.01401033338827501, 0.01401033338827501)
)

$startRow = 2
for ($i = 0; $i -lt $rows.Length; $i++) {
    $rowValues = $rows[$i]
    $targetRow = $startRow + $i
    for ($col = 0; $col -lt $rowValues.Length; $col++) {
        $ws.Cells.Item($targetRow, $col + 1).Value = $rowValues[$col]
    }
}

Write-Host "Updated $($rows.Length) data rows (A2:T9)."
